# Textbox response formatting fix
# Renames task-order sheets and updates the stimulus-file names referenced
# in column B of each sheet to a newer timestamped batch.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (new timestamped task-order identifiers) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687955517156"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687978006673"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168797802669"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687978480496"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687979225714"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687955184298.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687955356123.csv"
$ws1.Range("B4").Value = "go_stims-16511687955366578.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687955513039.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_2-1651168796010915.csv"
$ws2.Range("B3").Value = "OB-16511687969974542.csv"
$ws2.Range("B4").Value = "ZB-match_7-1651168796438322.csv"
$ws2.Range("B5").Value = "OB-16511687967640924.csv"
$ws2.Range("B6").Value = "TB-16511687977768023.csv"
$ws2.Range("B7").Value = "ZB-match_9-16511687958744051.csv"
$ws2.Range("B8").Value = "TB-16511687977093174.csv"
$ws2.Range("B9").Value = "OB-16511687972856417.csv"
$ws2.Range("B10").Value = "TB-16511687973971684.csv"

# --- Sheet 3: RS_TO --- (no cell content changes, only the sheet name above)

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687978160298.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687978046694.csv"
$ws4.Range("B4").Value = "MM_stims-16511687978313417.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687978170109.csv"
$ws4.Range("B6").Value = "MM_stims-16511687978469663.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687978323357.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687978782344.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168797907819.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687978528929.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687978920105.csv"
